$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: replace text dates with real Excel date serials, formatted as dates ---
$ws.Range("A2").Value = 45874
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy()
$ws.Range("A3:A7").PasteSpecial(-4122)

$ws.Range("A3").Value = 45875
$ws.Range("A4").Value = 45876
$ws.Range("A5").Value = 45876
$ws.Range("A6").Value = 45876
$ws.Range("A7").Value = 45876

# --- Column B: updated/reclassified description text ---
$ws.Range("B2").Value = "INTERMEDICA"
$ws.Range("B3").Value = "ASHS INTERMEDICA ASA"
$ws.Range("B4").Value = "Amil"
$ws.Range("B5").Value = "Unimed"
$ws.Range("B6").Value = "Sulamerica"
$ws.Range("B7").Value = "Bradesco"

# --- Column C: updated values ---
$ws.Range("C2").Value = 8
$ws.Range("C3").Value = 600
$ws.Range("C4").Value = 67
$ws.Range("C5").Value = 500
$ws.Range("C6").Value = 2600
$ws.Range("C7").Value = 568

# --- Widen column A so the date values are fully visible ---
$ws.Columns.Item(1).ColumnWidth = 16.3

# --- Leave the selection where the workbook was last saved ---
$ws.Range("D10").Select()
